$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 347 was just a placeholder label for station "FMWT 520" with no
# coordinates yet. Delete it entirely -- this shifts every row below it
# (348-394) up by one (347-393), which is exactly the bulk of the diff.
$ws.Rows.Item(347).Delete()

# Give the new last data row (385) the format of the other "FMWT" rows in
# column A (style 32) ...
$ws.Range("A384").Copy()
$ws.Range("A385").PasteSpecial(-4122)
$ws.Range("A385").Value = "FMWT"

# ... and now populate the newly-geocoded station 520 with real coordinates,
# reusing the header row's number formats (style 28 for B, 29 for C/D).
$ws.Range("B385").NumberFormat = "@"
$ws.Range("B385").Value = "520"
$ws.Range("C385").Value = 38.0328055556
$ws.Range("D385").Value = -121.869305556

$ws.Range("B1:D1").Copy()
$ws.Range("B385:D385").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Match the saved file's selection state (whole row 347 selected).
$ws.Rows.Item(347).Select()
